# Insert a new data row at row 150 (pushing existing rows 150-234 down to
# 151-235), then populate the new row with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(150).Insert()

$ws.Cells.Item(150, 1).Value  = 11
$ws.Cells.Item(150, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(150, 3).Value  = "Bíobío"
$ws.Cells.Item(150, 4).Value  = 45086
$ws.Cells.Item(150, 5).Value  = 8
$ws.Cells.Item(150, 6).Value  = 100112032
$ws.Cells.Item(150, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(150, 8).Value  = "Sin especificar"
$ws.Cells.Item(150, 9).Value  = "Primera"
$ws.Cells.Item(150, 10).Value = 170
$ws.Cells.Item(150, 11).Value = 8000
$ws.Cells.Item(150, 12).Value = 9000
$ws.Cells.Item(150, 13).Value = 8529
$ws.Cells.Item(150, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(150, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(150, 16).Value = 171
$ws.Cells.Item(150, 17).Value = 50
$ws.Cells.Item(150, 18).Value = "Hortaliza"
